$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine last used data row (existing data occupies rows 2..52)
$lastRow = $ws.UsedRange.Rows.Count

# Header row: add new columns AD, AE, AF right after the existing AC column
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the style of the existing header cells (bold, centered, bordered)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the team's W/L/T record for every data row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 74   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 88   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
